$wb = $excel.ActiveWorkbook

# --- Rename "Test" sheet to "Seconds" ---
$secondsSheet = $wb.Worksheets.Item("Test")
$secondsSheet.Name = "Seconds"

# --- Add new "Days" sheet right after "Seconds" ---
$daysSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $secondsSheet)
$daysSheet.Name = "Days"

# Header cell - reuse existing shared string text so it maps to the same
# shared-string entry as the other sheets ("MUMPS Date/Time Instant").
$daysSheet.Range("A1").Value = "MUMPS Date/Time Instant"
$daysSheet.Range("A1").Font.Bold = $true
$daysSheet.Range("A1").NumberFormat = "0"

# Data rows 2-10: convert the "Answer" sheet's MUMPS seconds value into days
# by dividing by 86400 (seconds per day).
for ($r = 2; $r -le 10; $r++) {
    $cell = $daysSheet.Range("A$r")
    $cell.Formula = "=Answer!A$r/86400"
    $cell.NumberFormat = "0"
}

# Column width to roughly match the "Answer" sheet's first column.
$daysSheet.Columns("A").ColumnWidth = 23.8776

# Restore the selection that was on the "Days" sheet when it was saved.
$daysSheet.Range("C9").Select() | Out-Null
